# Generate Report for Handoff
# Update the localization status report: "b.md" is now ready for handoff
# (instead of "handed back: in sync with en-US"), with new handoff file
# references / timestamps, and an error detail noting the handback file
# version mismatch.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: row 3 (b.md) ---
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("E3").Value = "Ready for handoff"
$ws.Range("F3").Value = "Ready for handoff"
$ws.Range("G3").Value = "2016-08-26 20:36:49"

# --- zh-cn sheet: row 3 (b.md) ---
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("C3").Value = "Ready for handoff"
# "False" typed directly gets auto-coerced to a real Boolean by Excel; copy
# the existing text cell above (F2, already literal text "False") instead so
# the value round-trips as a shared string like the original file.
$ws.Range("F2").Copy($ws.Range("F3"))
$ws.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$ws.Range("H3").Value = "2016-08-26 20:36:45"
$ws.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6901a4226f0a2e0e4e37f5fab096fe9758bcdd7c/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f1ddbab1e4c8e1c4c73ba1670acd137fe11e186d/e2e/b.md."
$ws.Columns.Item(16).ColumnWidth = 40

# --- de-de sheet: row 3 (b.md) ---
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("F2").Copy($ws.Range("F3"))
$ws.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$ws.Range("H3").Value = "2016-08-26 20:36:49"
$ws.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6901a4226f0a2e0e4e37f5fab096fe9758bcdd7c/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f1ddbab1e4c8e1c4c73ba1670acd137fe11e186d/e2e/b.md."
$ws.Columns.Item(16).ColumnWidth = 40
